$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("E_I and A_N ratios")

# New row 24 duplicates the original row 23 (same mouse, layer, cell types,
# area flags) but keeps the original pair/cell number (2) and original
# comment text that used to live on row 23.
$ws.Range("A24").Value = $ws.Range("A23").Value2
$ws.Range("B24").Value = 2
$ws.Range("C24").Value = $ws.Range("C23").Value2
$ws.Range("D24").Value = $ws.Range("D23").Value2
$ws.Range("E24").Value = $ws.Range("E23").Value2
$ws.Range("F24").Value = $ws.Range("F23").Value2
$ws.Range("G24").Value = $ws.Range("G23").Value2
$ws.Range("H24").Value = $ws.Range("H23").Value2
$ws.Range("I24").Value = $ws.Range("I23").Value2

# Row 23 now has a pair/cell number of 1 and a new comment replacing the old
# one (which now lives on row 24).
$ws.Range("B23").Value = 1
$ws.Range("I23").Value = "Assess is poor."

$ws.Range("A25").Select()
